# Add hours for Monday and Tuesday on the "Week 1" sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Week 1")
$ws.Activate()

# Monday (row 4) hours for Georgia Fox (column D)
$ws.Range("D4").Value = 2
# Tuesday (row 5) hours for Georgia Fox (column D)
$ws.Range("D5").Value = 4.5

# Select column D (mirrors the click on the column header seen in the diff)
$ws.Range("D1:D1048576").Select()

$wb.Save()
